# Refresh the cryptos price/volume snapshot to match the latest coinranking.com pull.
# (Two rows - Polkadot/BinanceUSD, HuobiToken/Filecoin, Decentraland/WEMIXTOKEN - swapped order too.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.775.56'
$ws.Range("E2").Value = '  -1.27%  '

# Row 3
$ws.Range("D3").Value = '1.792.93'
$ws.Range("E3").Value = '  +0.02%  '

# Row 4
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  +0.44%  '

# Row 5
$ws.Range("D5").Value = '''1.005'
$ws.Range("E5").Value = '  +0.33%  '

# Row 6
$ws.Range("D6").Value = '''303.97'
$ws.Range("E6").Value = '  -3.15%  '

# Row 7
$ws.Range("D7").Value = '''0.4951'
$ws.Range("E7").Value = '  -4.54%  '

# Row 8
$ws.Range("D8").Value = '''0.3831'
$ws.Range("E8").Value = '  +0.63%  '

# Row 9
$ws.Range("D9").Value = '''0.09213'
$ws.Range("E9").Value = '  +15.11%  '

# Row 10
$ws.Range("D10").Value = '''1.088'
$ws.Range("E10").Value = '  -0.52%  '

# Row 11
$ws.Range("D11").Value = '''40.46'
$ws.Range("E11").Value = '  -2.33%  '

# Row 12
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = '''1.007'
$ws.Range("E12").Value = '  +0.38%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''6.242'
$ws.Range("E13").Value = '  -0.61%  '

# Row 14
$ws.Range("D14").Value = '''20.54'
$ws.Range("E14").Value = '  +0.31%  '

# Row 15
$ws.Range("D15").Value = '1.803.44'
$ws.Range("E15").Value = '  +1.08%  '

# Row 16
$ws.Range("D16").Value = '''7.143'
$ws.Range("E16").Value = '  -1.97%  '

# Row 17
$ws.Range("D17").Value = '''91.95'
$ws.Range("E17").Value = '  +0.44%  '

# Row 18
$ws.Range("D18").Value = '''0.00001099'
$ws.Range("E18").Value = '  +0.66%  '

# Row 19
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("E20").Value = '  +0.30%  '

# Row 21
$ws.Range("D21").Value = '''16.99'
$ws.Range("E21").Value = '  -1.82%  '

# Row 22
$ws.Range("D22").Value = '''5.922'
$ws.Range("E22").Value = '  -0.45%  '

# Row 23
$ws.Range("D23").Value = '27.873.86'
$ws.Range("E23").Value = '  -1.05%  '

# Row 24
$ws.Range("D24").Value = '''10.92'
$ws.Range("E24").Value = '  -1.87%  '

# Row 25
$ws.Range("D25").Value = '''2.225'
$ws.Range("E25").Value = '  -2.12%  '

# Row 26
$ws.Range("D26").Value = '''158.35'
$ws.Range("E26").Value = '  -1.28%  '

# Row 27
$ws.Range("D27").Value = '2.006.21'
$ws.Range("E27").Value = '  +0.65%  '

# Row 28
$ws.Range("D28").Value = '''20.33'
$ws.Range("E28").Value = '  -0.50%  '

# Row 29
$ws.Range("D29").Value = '''2.362'
$ws.Range("E29").Value = '  +1.48%  '

# Row 30
$ws.Range("D30").Value = '''126.29'
$ws.Range("E30").Value = '  +2.78%  '

# Row 31
$ws.Range("D31").Value = '''0.1067'
$ws.Range("E31").Value = '  -0.94%  '

# Row 32
$ws.Range("D32").Value = '''1.042'
$ws.Range("E32").Value = '  -1.11%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.514'
$ws.Range("E33").Value = '  -0.47%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''3.610'
$ws.Range("E34").Value = '  -1.65%  '

# Row 35
$ws.Range("D35").Value = '''0.06785'
$ws.Range("E35").Value = '  -5.51%  '

# Row 36
$ws.Range("D36").Value = '''8.796'
$ws.Range("E36").Value = '  +1.91%  '

# Row 37
$ws.Range("D37").Value = '''0.02288'
$ws.Range("E37").Value = '  -0.77%  '

# Row 38
$ws.Range("D38").Value = '''0.2124'
$ws.Range("E38").Value = '  -0.98%  '

# Row 39
$ws.Range("D39").Value = '''11.28'
$ws.Range("E39").Value = '  -6.04%  '

# Row 40
$ws.Range("D40").Value = '''4.904'
$ws.Range("E40").Value = '  -3.05%  '

# Row 41
$ws.Range("D41").Value = '''0.6090'
$ws.Range("E41").Value = '  -1.15%  '

# Row 42
$ws.Range("D42").Value = '''1.004'
$ws.Range("E42").Value = '  +0.36%  '

# Row 43
$ws.Range("D43").Value = '''1.138'
$ws.Range("E43").Value = '  -2.18%  '

# Row 44
$ws.Range("D44").Value = '''12.96'
$ws.Range("E44").Value = '  -1.88%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '''1.279'
$ws.Range("E45").Value = '  -2.68%  '

# Row 46
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.5813'
$ws.Range("E46").Value = '  -2.44%  '

# Row 47
$ws.Range("D47").Value = '''3.649'
$ws.Range("E47").Value = '  -3.01%  '

# Row 48
$ws.Range("D48").Value = '''123.16'
$ws.Range("E48").Value = '  -3.38%  '

# Row 49
$ws.Range("D49").Value = '''1.929'
$ws.Range("E49").Value = '  +0.56%  '

# Row 50
$ws.Range("D50").Value = '''1.165'
$ws.Range("E50").Value = '  -4.01%  '

# Row 51
$ws.Range("D51").Value = '''0.06713'
$ws.Range("E51").Value = '  -0.49%  '
